$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 82
$ws.Range("C2").Value = 38
$ws.Range("E2").Value = 68.33333333333333
$ws.Range("F2").Value = 23.870164
$ws.Range("G2").Value = 0.44719
$ws.Range("H2").Value = 0.04938387414541062
$ws.Range("I2").Value = 0.09679239332500482
$ws.Range("J2").Value = 23.966956393325
$ws.Range("K2").Value = 23.773371606675
